$d = $word.ActiveDocument
$d.Content.Find.Execute("Lab section: _______", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Lab section: ___L2A_", 2)
